# The workbook gained one new data row for "Ajo" (Garlic) prices at the
# Mercado Mayorista Lo Valledor de Santiago market. The new record was
# inserted right before the existing row that is currently row 256
# (date 2021-06-11 / serial 44358), pushing that row and every row below
# it down by one. As a consequence the sheet's used range grows from
# A1:R362 to A1:R363.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 256, shifting rows 256-362
# down to 257-363 (xlShiftDown = -4121).
$ws.Rows.Item(256).Insert(-4121)

# Populate the newly inserted row 256 with the new price record.
$ws.Range("A256").Value = 6
$ws.Range("B256").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C256").Value = "Metropolitana"
$ws.Range("D256").Value = 44489
$ws.Range("E256").Value = 13
$ws.Range("F256").Value = 100112003
$ws.Range("G256").Value = "Ajo"
$ws.Range("H256").Value = "Chino"
$ws.Range("I256").Value = "Primera"
$ws.Range("J256").Value = 1700
$ws.Range("K256").Value = 17000
$ws.Range("L256").Value = 18000
$ws.Range("M256").Value = 17294
$ws.Range("N256").Value = "$/caja 10 kilos"
$ws.Range("O256").Value = "China"
$ws.Range("P256").Value = 1729
$ws.Range("Q256").Value = 10
$ws.Range("R256").Value = "Hortaliza"
